# lesson-20.xlsx word-list reorder
#
# The vocabulary rows in the "Sheet" worksheet were shuffled into a new
# study order. Columns A (English) and B (Japanese) move together, row by
# row, within two regions of the sheet:
#
#   rows  2-41  : old sub-blocks [20-30][9-19][31-41][2-8]      -> new order
#   rows 55-117 : old sub-blocks [84-95][108-117][55-68][96-107][69-83] -> new order
#
# Row 1 (header) and rows 42-54 keep their original content.
#
# Because this is a pure permutation of existing rows, every source value
# is snapshotted first (via Range.Text, which this host resolves reliably
# as a plain string) before any cell gets overwritten, then the snapshot
# is written back out in the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$snapA = @{}
$snapB = @{}

foreach ($r in 2..41) {
    $snapA[$r] = $ws.Cells.Item($r, 1).Text
    $snapB[$r] = $ws.Cells.Item($r, 2).Text
}

foreach ($r in 55..117) {
    $snapA[$r] = $ws.Cells.Item($r, 1).Text
    $snapB[$r] = $ws.Cells.Item($r, 2).Text
}

# Region 1: destination rows 2-41 pull from old sub-blocks, in this order:
#   20-30, 9-19, 31-41, 2-8
$region1 = @(
    ,(20, 30)
    ,(9, 19)
    ,(31, 41)
    ,(2, 8)
)

$dest = 2
foreach ($block in $region1) {
    $lo = $block[0]
    $hi = $block[1]
    foreach ($src in $lo..$hi) {
        $ws.Cells.Item($dest, 1).Value = $snapA[$src]
        $ws.Cells.Item($dest, 2).Value = $snapB[$src]
        $dest = $dest + 1
    }
}

# Region 2: destination rows 55-117 pull from old sub-blocks, in this order:
#   84-95, 108-117, 55-68, 96-107, 69-83
$region2 = @(
    ,(84, 95)
    ,(108, 117)
    ,(55, 68)
    ,(96, 107)
    ,(69, 83)
)

$dest = 55
foreach ($block in $region2) {
    $lo = $block[0]
    $hi = $block[1]
    foreach ($src in $lo..$hi) {
        $ws.Cells.Item($dest, 1).Value = $snapA[$src]
        $ws.Cells.Item($dest, 2).Value = $snapB[$src]
        $dest = $dest + 1
    }
}
